# Update results from parameter search
# - Rows 61 and 66 were re-run with num_seeds=10 (previously 5 and 7) and get
#   refreshed metric values.
# - Two parameter combinations that were missing from the grid
#   (folder_unf=1/folder_lr=0.0005/folder_dropout=0.25 and
#    folder_unf=2/folder_lr=0.0005/folder_dropout=0.25) are inserted into
#   their sorted positions (new rows 70 and 79), pushing subsequent rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ResultRow($rowNum, $values) {
    $ws.Cells.Item($rowNum, 1).Value = "SGD"
    for ($i = 1; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($rowNum, $i + 1).Value = $values[$i]
    }
}

# Insert the two new rows first, in ascending order, so later inserts use
# final target row numbers.
$ws.Rows("70:70").Insert()
$ws.Rows("79:79").Insert()

# Row 61: folder_unf=0, folder_lr=0.0005, folder_dropout=0.25 (num_seeds 5 -> 10)
Set-ResultRow 61 @(25, 0, 0.0005, 0.25, 10, 0.41859, 0.01537, 1314.836, 504.03407, 291.926, 113.06308, 0.87157, 0.01046, 185.7, 5.37587, 41.3, 5.37587, 21.5, 4.64878, 240.5, 4.64878, 0.88451, 0.009310000000000001)

# Row 66: folder_unf=1, folder_lr=0.0001, folder_dropout=0.2 (num_seeds 7 -> 10)
Set-ResultRow 66 @(25, 1, 0.0001, 0.2, 10, 0.44274, 0.05049, 1564.504, 549.01656, 314.439, 111.75621, 0.85031, 0.01215, 180.4, 5.96657, 46.6, 5.96657, 26.6, 6.6366, 235.4, 6.6366, 0.86537, 0.01157)

# New row 70: folder_unf=1, folder_lr=0.0005, folder_dropout=0.25
Set-ResultRow 70 @(25, 1, 0.0005, 0.25, 10, 0.20612, 0.03054, 1215.185, 575.73763, 243.807, 116.91754, 0.9222900000000001, 0.01268, 204.4, 5.85377, 22.6, 5.85377, 15.4, 3.20416, 246.6, 3.20416, 0.92853, 0.01115)

# New row 79: folder_unf=2, folder_lr=0.0005, folder_dropout=0.25
Set-ResultRow 79 @(25, 2, 0.0005, 0.25, 10, 0.18661, 0.02877, 1144.936, 751.55887, 208.161, 138.08879, 0.93333, 0.01007, 206.4, 4.85798, 20.6, 4.85798, 12, 2.94392, 250, 2.94392, 0.9388300000000001, 0.00894)
